$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 435, shifting existing data
# (rows 435-456) down to rows 437-458.
$ws.Range("A435:A436").EntireRow.Insert()

# New row 435: Camote, 1a (cosecha), date 2022-02-18 (serial 44610)
$ws.Range("A435").Value = 3
$ws.Range("B435").Value = "Femacal de La Calera"
$ws.Range("C435").Value = "Coquimbo"
$ws.Range("D435").Value = 44610
$ws.Range("E435").Value = 5
$ws.Range("F435").Value = 100112045
$ws.Range("G435").Value = "Zapallo"
$ws.Range("H435").Value = "Camote"
$ws.Range("I435").Value = "1a (cosecha)"
$ws.Range("J435").Value = 130
$ws.Range("K435").Value = 400
$ws.Range("L435").Value = 430
$ws.Range("M435").Value = 414
$ws.Range("N435").Value = '$/kilo (volumen en unidades)'
$ws.Range("O435").Value = "Provincia de Talca"
$ws.Range("P435").Value = 414
$ws.Range("Q435").Value = 1
$ws.Range("R435").Value = "Hortaliza"

# New row 436: Paine, 1a (cosecha), same date 2022-02-18 (serial 44610)
$ws.Range("A436").Value = 3
$ws.Range("B436").Value = "Femacal de La Calera"
$ws.Range("C436").Value = "Coquimbo"
$ws.Range("D436").Value = 44610
$ws.Range("E436").Value = 5
$ws.Range("F436").Value = 100112045
$ws.Range("G436").Value = "Zapallo"
$ws.Range("H436").Value = "Paine"
$ws.Range("I436").Value = "1a (cosecha)"
$ws.Range("J436").Value = 80
$ws.Range("K436").Value = 230
$ws.Range("L436").Value = 230
$ws.Range("M436").Value = 230
$ws.Range("N436").Value = '$/kilo (volumen en unidades)'
$ws.Range("O436").Value = "Provincia de Talca"
$ws.Range("P436").Value = 230
$ws.Range("Q436").Value = 1
$ws.Range("R436").Value = "Hortaliza"
